# "Generate Report for Archive"
# - Status of the 789ed9b5-... file moves from "Ready for handoff" to
#   "In Translation" on the Overview sheet (columns E/F, one per locale)
#   and on each per-locale sheet (zh-cn, de-de) in the Status column (C).
# - The Status column on each sheet narrows to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update every cell that held the old "Ready for handoff" status so the
# shared string itself is replaced everywhere it appears.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Re-fit the now-narrower Status columns.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
